$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 4237.6
$ws.Cells.Item(32, 9).Value = 3532.6667
$ws.Cells.Item(32, 10).Value = 5295
$ws.Cells.Item(32, 11).Value = 3532.6667
$ws.Cells.Item(32, 12).Value = 5295
$ws.Cells.Item(32, 13).Value = -3206.6667
$ws.Cells.Item(32, 14).Value = -5947
# Row 53
$ws.Cells.Item(53, 8).Value = 329.6875
$ws.Cells.Item(53, 9).Value = 269.57144
$ws.Cells.Item(53, 10).Value = 376.44446
$ws.Cells.Item(53, 11).Value = 269.57144
$ws.Cells.Item(53, 12).Value = 376.44446
$ws.Cells.Item(53, 13).Value = 367.42856
$ws.Cells.Item(53, 14).Value = -1650.44446
# Row 76
$ws.Cells.Item(76, 8).Value = 4912.1665
$ws.Cells.Item(76, 10).Value = 4000
$ws.Cells.Item(76, 12).Value = 4000
$ws.Cells.Item(76, 14).Value = -4630
# Row 79
$ws.Cells.Item(79, 8).Value = 4912.1665
$ws.Cells.Item(79, 10).Value = 4000
$ws.Cells.Item(79, 12).Value = 4000
$ws.Cells.Item(79, 14).Value = -6184
# Row 112
$ws.Cells.Item(112, 8).Value = 2771.4443
$ws.Cells.Item(112, 10).Value = 3034.7144
$ws.Cells.Item(112, 12).Value = 9104.143199999999
$ws.Cells.Item(112, 14).Value = -11320.1432
# Row 132
$ws.Cells.Item(132, 8).Value = 2486.6667
$ws.Cells.Item(132, 9).Value = 1675.1515
$ws.Cells.Item(132, 10).Value = 6950
$ws.Cells.Item(132, 11).Value = 5025.4545
$ws.Cells.Item(132, 12).Value = 20850
$ws.Cells.Item(132, 13).Value = -2495.4545
$ws.Cells.Item(132, 14).Value = -25910
# Row 138
$ws.Cells.Item(138, 8).Value = 2596.3
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2300.5715
$ws.Cells.Item(32, 9).Value = 2300.5715
$ws.Cells.Item(32, 11).Value = 2300.5715
$ws.Cells.Item(32, 13).Value = -2013.5715
# Row 61
$ws.Cells.Item(61, 8).Value = 2836.077
$ws.Cells.Item(61, 9).Value = 2836.077
$ws.Cells.Item(61, 11).Value = 2836.077
$ws.Cells.Item(61, 13).Value = -2624.077
# Row 122
$ws.Cells.Item(122, 8).Value = 2430.5
$ws.Cells.Item(122, 9).Value = 2273.6667
$ws.Cells.Item(122, 11).Value = 6821.000100000001
$ws.Cells.Item(122, 13).Value = -4371.000100000001
# Row 132
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
# Row 136
$ws.Cells.Item(136, 8).Value = 2836.077
$ws.Cells.Item(136, 9).Value = 2836.077
$ws.Cells.Item(136, 11).Value = 8508.231
$ws.Cells.Item(136, 13).Value = -5958.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 2114.7144
$ws.Cells.Item(20, 9).Value = 2114.7144
$ws.Cells.Item(20, 11).Value = 2114.7144
$ws.Cells.Item(20, 13).Value = -1867.7144
# Row 95
$ws.Cells.Item(95, 8).Value = 18312.5
$ws.Cells.Item(95, 10).Value = 18312.5
$ws.Cells.Item(95, 12).Value = 18312.5
$ws.Cells.Item(95, 14).Value = -23804.5
# Row 99
$ws.Cells.Item(99, 8).Value = 1250.2222
$ws.Cells.Item(99, 9).Value = 1040.1333
$ws.Cells.Item(99, 11).Value = 1040.1333
$ws.Cells.Item(99, 13).Value = 457.8667
# Row 107
$ws.Cells.Item(107, 8).Value = 5098
$ws.Cells.Item(107, 9).Value = 3684.875
$ws.Cells.Item(107, 11).Value = 3684.875
$ws.Cells.Item(107, 13).Value = -1764.875
# Row 134
$ws.Cells.Item(134, 8).Value = 2898.577
$ws.Cells.Item(134, 9).Value = 2898.577
$ws.Cells.Item(134, 11).Value = 8695.731
$ws.Cells.Item(134, 13).Value = -6160.731

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2461.2778
$ws.Cells.Item(16, 10).Value = 3779.2
$ws.Cells.Item(16, 12).Value = 3779.2
$ws.Cells.Item(16, 14).Value = -4353.2
# Row 94
$ws.Cells.Item(94, 8).Value = 147114.38
$ws.Cells.Item(94, 9).Value = 228981.6
$ws.Cells.Item(94, 10).Value = 10669
$ws.Cells.Item(94, 11).Value = 228981.6
$ws.Cells.Item(94, 12).Value = 10669
$ws.Cells.Item(94, 13).Value = -228530.6
$ws.Cells.Item(94, 14).Value = -11571
# Row 113
$ws.Cells.Item(113, 8).Value = 2461.2778
$ws.Cells.Item(113, 10).Value = 3779.2
$ws.Cells.Item(113, 12).Value = 3779.2
$ws.Cells.Item(113, 14).Value = -8119.2
# Row 122
$ws.Cells.Item(122, 8).Value = 955.5
$ws.Cells.Item(122, 9).Value = 955.5
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2866.5
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -416.5
$ws.Cells.Item(122, 14).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 1752.04
$ws.Cells.Item(132, 9).Value = 1811.2778
$ws.Cells.Item(132, 10).Value = 1599.7142
$ws.Cells.Item(132, 11).Value = 5433.8334
$ws.Cells.Item(132, 12).Value = 4799.142599999999
$ws.Cells.Item(132, 13).Value = -2903.8334
$ws.Cells.Item(132, 14).Value = -9859.142599999999
# Row 134
$ws.Cells.Item(134, 8).Value = 721.2759
$ws.Cells.Item(134, 9).Value = 737.75
$ws.Cells.Item(134, 11).Value = 2213.25
$ws.Cells.Item(134, 13).Value = 321.75
# Row 141
$ws.Cells.Item(141, 8).Value = 35677.875
$ws.Cells.Item(141, 10).Value = 35057.133
$ws.Cells.Item(141, 12).Value = 35057.133
$ws.Cells.Item(141, 14).Value = -45417.133

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 362
$ws.Cells.Item(7, 9).Value = 211.25
$ws.Cells.Item(7, 10).Value = 462.5
$ws.Cells.Item(7, 11).Value = 633.75
$ws.Cells.Item(7, 12).Value = 1387.5
$ws.Cells.Item(7, 13).Value = -521.75
$ws.Cells.Item(7, 14).Value = -1611.5
# Row 98
$ws.Cells.Item(98, 8).Value = 2500
$ws.Cells.Item(98, 10).Value = 2500
$ws.Cells.Item(98, 12).Value = 7500
$ws.Cells.Item(98, 14).Value = -10496
# Row 121
$ws.Cells.Item(121, 8).Value = 12384.571
$ws.Cells.Item(121, 9).Value = 38686.668
$ws.Cells.Item(121, 10).Value = 5211.273
$ws.Cells.Item(121, 11).Value = 116060.004
$ws.Cells.Item(121, 12).Value = 15633.819
$ws.Cells.Item(121, 13).Value = -114750.004
$ws.Cells.Item(121, 14).Value = -18253.819

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 9476.5
$ws.Cells.Item(70, 9).Value = 8976
$ws.Cells.Item(70, 10).Value = 9977
$ws.Cells.Item(70, 11).Value = 8976
$ws.Cells.Item(70, 12).Value = 9977
$ws.Cells.Item(70, 13).Value = -8706
$ws.Cells.Item(70, 14).Value = -10517
# Row 73
$ws.Cells.Item(73, 8).Value = 9476.5
$ws.Cells.Item(73, 9).Value = 8976
$ws.Cells.Item(73, 10).Value = 9977
$ws.Cells.Item(73, 11).Value = 8976
$ws.Cells.Item(73, 12).Value = 9977
$ws.Cells.Item(73, 13).Value = -8040
$ws.Cells.Item(73, 14).Value = -11849
# Row 107
$ws.Cells.Item(107, 8).Value = 1982.5
$ws.Cells.Item(107, 10).Value = 2810.3333
$ws.Cells.Item(107, 12).Value = 2810.3333
$ws.Cells.Item(107, 14).Value = -6650.3333
# Row 122
$ws.Cells.Item(122, 8).Value = 3136.9473
$ws.Cells.Item(122, 9).Value = 3133.5334
$ws.Cells.Item(122, 10).Value = 3149.75
$ws.Cells.Item(122, 11).Value = 9400.600199999999
$ws.Cells.Item(122, 12).Value = 9449.25
$ws.Cells.Item(122, 13).Value = -6950.600199999999
$ws.Cells.Item(122, 14).Value = -14349.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5295.125
$ws.Cells.Item(7, 9).Value = 1622.5385
$ws.Cells.Item(7, 11).Value = 1622.5385
$ws.Cells.Item(7, 13).Value = -1510.5385
# Row 100
$ws.Cells.Item(100, 8).Value = 3026.4285
$ws.Cells.Item(100, 9).Value = 2784.0715
$ws.Cells.Item(100, 10).Value = 3511.1428
$ws.Cells.Item(100, 11).Value = 2784.0715
$ws.Cells.Item(100, 12).Value = 3511.1428
$ws.Cells.Item(100, 13).Value = -2243.0715
$ws.Cells.Item(100, 14).Value = -4593.1428
# Row 126
$ws.Cells.Item(126, 8).Value = 5295.125
$ws.Cells.Item(126, 9).Value = 1622.5385
$ws.Cells.Item(126, 11).Value = 4867.6155
$ws.Cells.Item(126, 13).Value = -2397.6155
# Row 136
$ws.Cells.Item(136, 8).Value = 4639.857
$ws.Cells.Item(136, 9).Value = 4639.857
$ws.Cells.Item(136, 11).Value = 13919.571
$ws.Cells.Item(136, 13).Value = -11369.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Cells.Item(135, 8).Value = 363404.34
$ws.Cells.Item(135, 10).Value = 363404.34
$ws.Cells.Item(135, 12).Value = 363404.34
$ws.Cells.Item(135, 14).Value = -373544.34
